$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1 becomes a text label instead of the numeric 0 it held before, while
# keeping its existing (bold/bordered header) style.
$ws.Range("B1").Value = "Money Raised Currency (in USD)"

# A2:A66 lose the bold/bordered header style that had been applied to them
# (only row 1 keeps that header formatting now).
$ws.Range("A2:A66").Style = "Normal"
